$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 52

$ws.Cells.Item($row, 1).Value2  = 11
$ws.Cells.Item($row, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value2  = "Bíobío"
$ws.Cells.Item($row, 4).Value2  = 45191
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
$ws.Cells.Item($row, 5).Value2  = 8
$ws.Cells.Item($row, 6).Value2  = "Fruta"
$ws.Cells.Item($row, 7).Value2  = 100107
$ws.Cells.Item($row, 8).Value2  = "Otros"
$ws.Cells.Item($row, 9).Value2  = 100107002
$ws.Cells.Item($row, 10).Value2 = "Chirimoya"
$ws.Cells.Item($row, 11).Value2 = "Cultivar IV Región"
$ws.Cells.Item($row, 12).Value2 = "Primera"
$ws.Cells.Item($row, 13).Value2 = 140
$ws.Cells.Item($row, 14).Value2 = 20000
$ws.Cells.Item($row, 15).Value2 = 21000
$ws.Cells.Item($row, 16).Value2 = 20571
$ws.Cells.Item($row, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item($row, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item($row, 19).Value2 = 2057
$ws.Cells.Item($row, 20).Value2 = 10
